# Regenerate the "K" column (column G) values on the active sheet.
# This mirrors the commit "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" -- the recalculated strikeout
# counts (K) replace the previous placeholder values for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, taken from the recalculated
# save_data output.
$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 3
    12 = 0
    13 = 1
    15 = 1
    17 = 3
    18 = 2
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
